$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Holly added "S.GISH" as a harvester value in bioSamples; fix the harvester
# column (B) in rnaSamples so every data row reflects the same harvester.
$ws.Range("B2:B16").Value = "S.GISH"

# Leave the selection on the harvester column, matching the edited state.
[void]$ws.Range("B:B").Select()
